$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("状态" / status) values for rows 2-55, added/updated per the
# "add some china word" commit (rows 3-55 previously had no B value; B2
# changes from 3 to 4).
$statusValues = @(
    4,  # B2  (row 2, was 3)
    1,  # B3  (row 3)
    4,  # B4
    4,  # B5
    4,  # B6
    4,  # B7
    1,  # B8
    4,  # B9
    4,  # B10
    4,  # B11
    4,  # B12
    4,  # B13
    1,  # B14
    4,  # B15
    4,  # B16
    1,  # B17
    1,  # B18
    4,  # B19
    4,  # B20
    4,  # B21
    4,  # B22
    4,  # B23
    1,  # B24
    1,  # B25
    4,  # B26
    4,  # B27
    4,  # B28
    4,  # B29
    4,  # B30
    4,  # B31
    4,  # B32
    4,  # B33
    4,  # B34
    4,  # B35
    4,  # B36
    4,  # B37
    4,  # B38
    4,  # B39
    4,  # B40
    1,  # B41
    4,  # B42
    4,  # B43
    4,  # B44
    4,  # B45
    4,  # B46
    4,  # B47
    4,  # B48
    4,  # B49
    4,  # B50
    1,  # B51
    1,  # B52
    1,  # B53
    1,  # B54
    1   # B55
)

$startRow = 2
for ($i = 0; $i -lt $statusValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $statusValues[$i]
}

# Reposition the view/selection like the saved workbook: scrolled down so
# row 38 is at the top, with D47 as the active selected cell.
$excel.Goto($ws.Range("A38"), $true)
$ws.Range("D47").Select()
